$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume/ranking data per the Dec 29 2023 refresh.
# Cells whose new value is a plain numeric-looking string are written with a
# leading apostrophe to force text (matching the original inline-string cells),
# then the style is reset to "Normal" so no stray number-format is introduced.

$ws.Range("D2").Value = "42.801.53"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.358.86"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'319.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "'106.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("E7").Value = "  -2.43%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -4.73%  "
$ws.Range("D10").Value = "'41.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "'8.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'15.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.80%  "
$ws.Range("D16").Value = "2.712.22"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "2.414.26"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "42.789.52"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "'7.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.57%  "
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "'76.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "'3.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.40%  "
$ws.Range("D23").Value = "'258.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.84%  "
$ws.Range("E24").Value = "  -4.54%  "
$ws.Range("D25").Value = "'9.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.41%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").Value = "'23.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("D30").Value = "'174.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'36.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("E32").Value = "  +3.94%  "
$ws.Range("D33").Value = "'0.0887"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.56%  "
$ws.Range("E34").Value = "  -8.00%  "
$ws.Range("E35").Value = "  +14.65%  "
$ws.Range("D36").Value = "'0.131"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").Value = "'4.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.08%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("E39").Value = "  -8.92%  "
$ws.Range("D40").Value = "'2.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.81%  "
$ws.Range("D41").Value = "'0.237"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").Value = "'71.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("E43").Value = "  -7.92%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'12.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.21%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'114.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.55%  "
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'85.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'9.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.17%  "
$ws.Range("D50").Value = "'74.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("E51").Value = "  -1.69%  "
